# Applies the "storage::manager" -> "storage::api" rename (and the related
# *Manager -> *Storage box relabels + connector nudge) to the Storage
# component diagram slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "storage::manager" -> "storage::api" -------------------------------
# This text lives on "Rectangle 12", the first item of the "Group 11"
# group shape. It is currently a single run; split "manager" off into its
# own run (by re-typing it in place) so the result has the three runs
# "storage" / "::" / "api" like the target XML.
$grp = $s.Shapes.Item(2)
$lbl = $grp.GroupItems.Item(1)
$tr = $lbl.TextFrame.TextRange
$mid = $tr.Characters(8, 2)
$mid.Text = "::"
$tail = $tr.Characters(10, 7)
$tail.Text = "api"

# --- nudge the connector that points at the storage::api box ------------
$conn = $s.Shapes.Item(4)
$conn.Left = 82.11035

# --- *Manager -> *Storage box relabels -----------------------------------
$s.Shapes.Item(6).TextFrame.TextRange.Text = "EvaluationsStorage"
$s.Shapes.Item(7).TextFrame.TextRange.Text = "AccountsStorage"
$s.Shapes.Item(10).TextFrame.TextRange.Text = "CoursesStorage"
